$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells A1/J1 (text unchanged, but force rewrite)
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

# --- Negative table (columns A-H): now only 5 data rows (3-7) ---
$negData = @(
    @("crude", 0.7941176470588235, 27, 27, 0, 1, $false, 7),
    @("fraud", 0.6388888888888888, 23, 23, 0, 1, $false, 13),
    @("crisis", 0.6164383561643836, 180, 180, 0, 1, $false, 112),
    @("panic", 0.1802325581395349, 93, 93, 0, 1, $false, 423),
    @("sc", 0.1587301587301587, 30, 30, 0, 1, $false, 159)
)
for ($i = 0; $i -lt $negData.Count; $i++) {
    $r = 3 + $i
    $row = $negData[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
}

# Clear old rows 8-10 in columns A-H (table shrank from 8 to 5 data rows)
$ws.Range("A8:H10").Clear()

# Seed rows 37-38 in column J with the same style used by the rest of the
# "name" column (bold/bordered/centered header style) by copying from J36,
# then the value loop below overwrites the copied value with the real one.
$ws.Range("J36").Copy($ws.Range("J37"))
$ws.Range("J36").Copy($ws.Range("J38"))

# --- Positive table (columns J-Q): now 36 data rows (3-38) ---
$posData = @(
    @("best", 0.9491525423728814, 56, 56, 1, 0, $false, 3),
    @("interesting", 0.9393939393939394, 31, 31, 1, 0, $false, 2),
    @("happy", 0.9230769230769231, 24, 24, 1, 0, $false, 2),
    @("strong", 0.9090909090909091, 20, 20, 1, 0, $false, 2),
    @("love", 0.8913043478260869, 41, 41, 1, 0, $false, 5),
    @("nice", 0.8888888888888888, 24, 24, 1, 0, $false, 3),
    @("great", 0.8660714285714286, 97, 97, 1, 0, $false, 15),
    @("special", 0.8333333333333334, 30, 30, 1, 0, $false, 6),
    @("thank", 0.828125, 106, 106, 1, 0, $false, 22),
    @("thanks", 0.8170731707317073, 67, 67, 1, 0, $false, 15),
    @("healthy", 0.7777777777777778, 21, 21, 1, 0, $false, 6),
    @("positive", 0.7758620689655172, 45, 45, 1, 0, $false, 13),
    @("free", 0.7666666666666667, 92, 92, 1, 0, $false, 28),
    @("safety", 0.7254901960784313, 37, 37, 1, 0, $false, 14),
    @("safe", 0.7253521126760564, 103, 103, 1, 0, $false, 39),
    @("support", 0.7169811320754716, 76, 76, 1, 0, $false, 30),
    @("friends", 0.7142857142857143, 20, 20, 1, 0, $false, 8),
    @("good", 0.69375, 111, 111, 1, 0, $false, 49),
    @("fresh", 0.6666666666666666, 32, 32, 1, 0, $false, 16),
    @("confidence", 0.6388888888888888, 23, 23, 1, 0, $false, 13),
    @("better", 0.6349206349206349, 40, 40, 1, 0, $false, 23),
    @("well", 0.6170212765957447, 58, 58, 1, 0, $false, 36),
    @("relief", 0.6, 30, 30, 1, 0, $false, 20),
    @("heroes", 0.5957446808510638, 28, 28, 1, 0, $false, 19),
    @("hand", 0.4986945169712794, 191, 191, 1, 0, $false, 192),
    @("important", 0.4888888888888889, 22, 22, 1, 0, $false, 23),
    @("like", 0.4705882352941176, 160, 160, 1, 0, $false, 180),
    @("care", 0.4269662921348314, 38, 38, 1, 0, $false, 51),
    @("help", 0.4203389830508474, 124, 124, 1, 0, $false, 171),
    @("hope", 0.4, 26, 26, 1, 0, $false, 39),
    @("protect", 0.3698630136986301, 27, 27, 1, 0, $false, 46),
    @("please", 0.3514644351464435, 84, 84, 1, 0, $false, 155),
    @("sure", 0.3125, 20, 20, 1, 0, $false, 44),
    @("increase", 0.3076923076923077, 24, 24, 1, 0, $false, 54),
    @("store", 0.02237136465324385, 20, 20, 1, 0, $false, 874),
    @("co", 0.006437077566784679, 20, 20, 1, 0, $false, 3087)
)
for ($i = 0; $i -lt $posData.Count; $i++) {
    $r = 3 + $i
    $row = $posData[$i]
    $ws.Range("J$r").Value = $row[0]
    $ws.Range("K$r").Value = $row[1]
    $ws.Range("L$r").Value = $row[2]
    $ws.Range("M$r").Value = $row[3]
    $ws.Range("N$r").Value = $row[4]
    $ws.Range("O$r").Value = $row[5]
    $ws.Range("P$r").Value = $row[6]
    $ws.Range("Q$r").Value = $row[7]
}


